# Commit: "added templates for event types (#259)"
#
# Adds new "event type" email/notification templates to the template_type
# sheet, one row per (language, event-type) combination, appended after the
# existing data (which ended at row 2597).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Languages, in the same order used for the other multi-language template
# blocks already present in the sheet (eng, fra, ara, hin, kan, tam, spa).
$langs = @("eng", "fra", "ara", "hin", "kan", "tam", "spa")

# New event-type template codes (column B) and their descriptions (column C).
$codes = @(
    "mosip.event.type.AUTHENTICATION_REQUEST",
    "mosip.event.type.SHARE_CRED_WITH_PARTNER",
    "mosip.event.type.DOWNLOAD_PERSONALIZED_CARD",
    "mosip.event.type.ORDER_PHYSICAL_CARD",
    "mosip.event.type.GET_MY_ID",
    "mosip.event.type.BOOK_AN_APPOINTMENT",
    "mosip.event.type.UPDATE_MY_UIN",
    "mosip.event.type.GENERATE_VID",
    "mosip.event.type.REVOKE_VID",
    "mosip.event.type.AUTH_TYPE_LOCK_UNLOCK",
    "mosip.event.type.VID_CARD_DOWNLOAD",
    "mosip.event.type.SEND_OTP",
    "mosip.event.type.VALIDATE_OTP",
    "mosip.event.type.DEFAULT"
)

$descrs = @(
    "Authentication Request event type",
    "Share Credential With Partner event type",
    "Download Personalized Card event type",
    "Order a Physical Card event type",
    "Get UIN Card event type",
    "Book An Appointment event type",
    "Update UIN Data event type",
    "Generate VID event type",
    "Revoke VID event type",
    "Secure My ID event type",
    "Download VID Card event type",
    "Send OTP event type",
    "Verify My Phone/Email event type",
    "Default event type"
)

# Column D (is_active) already holds the literal text "TRUE" (not the
# boolean value) throughout the sheet, e.g. cell D2. Assigning the string
# "TRUE" straight to .Value would be auto-converted to a real boolean, so
# instead copy the existing text cell and paste just its value into each
# new D cell - this preserves both the "TRUE" text and the D-column style.
$trueTextCell = $ws.Range("D2")

$row = 2598

foreach ($lang in $langs) {
    for ($i = 0; $i -lt $codes.Length; $i++) {
        $ws.Cells.Item($row, 1).Value = $lang
        $ws.Cells.Item($row, 2).Value = $codes[$i]
        $ws.Cells.Item($row, 3).Value = $descrs[$i]
        $trueTextCell.Copy()
        $ws.Cells.Item($row, 4).PasteSpecial(-4163)
        $row = $row + 1
    }
}

# Keep the active selection pointing just past the newly written data,
# matching the author's workbook view after the edit.
$null = $ws.Range("A" + $row).Select()
